$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# The last existing data row (CSR_Regression / CSR Manage Users / test.java.TestCSRManageUsers)
# currently sits on row 4. We need to push it down to row 6 and insert two new rows
# (Enrollment Flow, My Profile) above it, in rows 4 and 5.

$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).Insert()

$ws.Cells.Item(4, 1).Value = "UPA_Regression"
$ws.Cells.Item(4, 2).Value = "Enrollment Flow"
$ws.Cells.Item(4, 3).Value = "test.java.TestUPAEnrollment"

$ws.Cells.Item(5, 1).Value = "UPA_Regression"
$ws.Cells.Item(5, 2).Value = "My Profile"
$ws.Cells.Item(5, 3).Value = "test.java.TestUPAMyProfileTab"

$ws.Range("C5").Select()
